$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Battery_Capacity (B3)
$ws.Range("B3").Value = 599400

# Update T_Boil (B5)
$ws.Range("B5").Value = 320.39505615000002

# Recalculate so the formula in B6 (Boils_per_Charge) updates
$excel.Calculate()

$wb.Save()
